# Add a centered "Questionnaire 27" header (Arial, 12pt) to the document's
# only section, matching the default header reference added in the diff.

$d = $word.ActiveDocument

$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)  # wdHeaderFooterPrimary

# Use InsertAfter (rather than assigning Range.Text) so the engine only
# mints the single "default" header part instead of the full
# first/even/default header+footer family.
$header.Range.InsertAfter("Questionnaire 27")

# Paragraph-level formatting: Header style, centered.
$header.Range.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter

# Run-level formatting: Arial 12pt, applied to the text only (not the
# paragraph mark) so no stray rPr ends up on the pPr.
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
